# Add two new columns: I (I0) and J (IF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new header labels in I1/J1, matching existing header style (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-20: I column is always 1, J column mirrors the H column value
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($r, 8).Value2
}
